$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column A (date) cells from Excel serial-date numbers to plain
# YYYYMMDD integers, and drop the date number-format style that was
# previously applied to them (style index 2 -> back to default/no style).

$ws.Range("A2:A7").Value = 20180405
$ws.Range("A8:A13").Value = 20180412
$ws.Range("A14:A19").Value = 20180418
$ws.Range("A20:A25").Value = 20180426
$ws.Range("A26:A31").Value = 20180503
$ws.Range("A32:A37").Value = 20180510
$ws.Range("A38:A43").Value = 20180517
$ws.Range("A44:A49").Value = 20180524
$ws.Range("A50:A55").Value = 20180531
$ws.Range("A56:A61").Value = 20180607
$ws.Range("A62:A67").Value = 20180614
$ws.Range("A68:A73").Value = 20180621

# Clear the per-cell style (removes the custom date numFmt) so the cells
# revert to the workbook default style, matching the author's change that
# removed the now-unused numFmt/cellXfs entries.
$ws.Range("A2:A73").Style = "Normal"
